$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Numeric data for the new columns (rows 2-56)
$iVals = @(1,1,1,1,1,1,1,1,6,4,1,6,8,6,8,6,8,7,7,6,8,6,3,9,10,8,5,9,7,7,7,6,7,9,7,9,8,7,9,7,8,7,6,5,7,6,6,5,5,6,4,4,5,6,4)
$jVals = @(4,6,7,5,5,5,4,3,7,7,2,8,8,7,8,6,8,7,7,7,9,8,6,9,10,9,6,9,8,8,8,7,8,10,8,9,9,8,9,8,9,9,7,7,9,7,7,6,6,7,6,5,6,6,5)

for ($r = 0; $r -lt $iVals.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
